$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Insert a new row at row 14, shifting existing rows 14+ down.
$ws.Rows.Item(14).Insert()

$ws.Range("B14").Value = "[STM32F4xxhttps://github.com/grblHAL/STM32F4xx"
$ws.Range("A14").Value = "[Flexi-HAL](https://github.com/Expatria-Technologies/Flexi-HAL) "

$ws.Range("B13").Select()
